$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '69.168.02', '  -2.57%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.672.39', '  -4.11%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  +0.12%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '673.56', '  -4.73%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '159.78', '  -6.90%  '),
    @(7, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.671.17', '  -4.07%  '),
    @(8, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  +0.19%  '),
    @(9, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.479', '  -8.44%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.147', '  -9.11%  '),
    @(11, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '7.11', '  -4.25%  '),
    @(12, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.444', '  -2.98%  '),
    @(13, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000227', '  -10.72%  '),
    @(14, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '4.282.77', '  -4.28%  '),
    @(15, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '32.71', '  -10.43%  '),
    @(16, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.687.26', '  -2.75%  '),
    @(17, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '69.011.27', '  -2.81%  '),
    @(18, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.112', '  -2.18%  '),
    @(19, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '16.12', '  -7.08%  '),
    @(20, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.48', '  -10.24%  '),
    @(21, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '477.84', '  -3.61%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '9.76', '  -8.38%  '),
    @(23, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.658', '  -10.28%  '),
    @(24, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '77.76', '  -8.94%  '),
    @(25, 'WrappedeETH', 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth', '3.803.49', '  -4.48%  '),
    @(26, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.04%  '),
    @(27, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '11.40', '  -5.71%  '),
    @(28, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0000125', '  -13.58%  '),
    @(29, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '9.09', '  -14.24%  '),
    @(30, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '1.82', '  -12.67%  '),
    @(31, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '2.69', '  -13.07%  '),
    @(32, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.03', '  -8.70%  '),
    @(33, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '6.59', '  -11.02%  '),
    @(34, 'Binance-PegBSC-USD', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd', '0.997', '  -0.22%  '),
    @(35, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '26.45', '  -9.84%  '),
    @(36, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.163', '  -7.21%  '),
    @(37, 'RenzoRestakedETH', 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth', '3.635.99', '  -4.30%  '),
    @(38, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.41', '  -8.23%  '),
    @(39, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.00', '  +0.44%  '),
    @(40, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0914', '  -10.58%  '),
    @(41, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '1.00', '  -0.03%  '),
    @(42, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '2.18', '  -6.16%  '),
    @(43, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.997', '  -0.29%  '),
    @(44, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.942', '  -10.21%  '),
    @(45, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '159.46', '  -2.43%  '),
    @(46, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '47.86', '  -2.21%  '),
    @(47, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '2.84', '  -14.41%  '),
    @(48, 'ONDO', 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo', '1.30', '  -4.78%  '),
    @(49, 'FLOKI', 'https://coinranking.com/coin/fmHk13Rqw+floki-floki', '0.000271', '  -13.01%  '),
    @(50, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '7.88', '  -9.89%  '),
    @(51, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '375.37', '  -12.67%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
